# TRIMAZKON_address_list.xlsx - "dodelany manage bin u ip setting v6"
#
# - row 1 (527_Teijin) renamed to "(2)" and its IP note trimmed
# - a brand-new "527_Teijin(1)" row is inserted as row 2, with its own
#   (sparser) IP note
# - the old duplicate "527_Teijin (1)" row (previously row 3) is removed
# - the 529_Witte camera IP note is shortened
# - the last row's "527_Teijin" label gets its "(1)" suffix back

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ip_address_list")

# --- Row 1: 527_Teijin (1) -> 527_Teijin (2), trimmed IP note ---
$ws.Range("A1").Value = "527_Teijin (2)"
$ws.Range("D1").Value = "XG-X290" + "`n" + "OP:`t`t10.101.28.11"

# --- Insert a fresh row above the old row 2, becomes the new row 2 ---
$ws.Rows(2).Insert()
$ws.Range("A2").Value = "527_Teijin(1)"
$ws.Range("B2").Value = "10.101.28.176"
$ws.Range("C2").Value = "255.255.255.0"
$ws.Range("D2").Value = "XG-X2900:`t`t10.101.28" + "`n" + "OP:`t`t10.101.28.11"
# keep E as text ("1"), matching the rest of the column, instead of a number
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1"

# --- Remove the now-duplicated 527_Teijin (1) row (old row 3, shifted to row 4) ---
$ws.Rows(4).Delete()

# --- Row 6 (529_Witte): shorten the camera IP note ---
$ws.Range("D6").Value = "Kamera VS-S160MX :192.168.0.1"

# --- Row 10: 527_Teijin -> 527_Teijin (1) ---
$ws.Range("A10").Value = "527_Teijin (1)"
